$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$rowData = @(
    @(1.02, 1.017102856055299, 1.028391008553387, 1.018544776484532, 1.015459048434933, 1.030222706043753, 1.022318879346831, 1.031207962458893, 1.021390675308571, 1.018314144156688, 1.011535190582853),
    @(1.02, 1.017955109450439, 1.028823679859129, 1.019264920005833, 1.016951933063509, 1.030338368904945, 1.022807384840872, 1.031449572333836, 1.021916888089709, 1.019610290616794, 1.011698472467801),
    @(1.02, 1.01850679083764, 1.029103856349498, 1.019731482279518, 1.01791811974935, 1.030412093533878, 1.023123076464071, 1.031605386702157, 1.022257296765029, 1.020448688472642, 1.01180394401526),
    @(1.02, 1.018738768946457, 1.029221690989423, 1.019927763266207, 1.018324354428606, 1.030442819504783, 1.023255696111768, 1.031670764801687, 1.022400383378366, 1.020801083677376, 1.01184824032295),
    @(1.02, 1.018777722080773, 1.029241478736703, 1.019960727784414, 1.018392566060602, 1.03044796280959, 1.023277957828781, 1.031681734640179, 1.022424406991152, 1.020860248540983, 1.011855675286541),
    @(1.02, 1.018509890339801, 1.029105430673461, 1.019734104453811, 1.017923547676128, 1.030412505148895, 1.023124848916342, 1.03160626078486, 1.022259208779349, 1.02045339745451, 1.011804536077538),
    @(1.02, 1.017390833923287, 1.028537187539103, 1.018788031436844, 1.015963538861435, 1.030262025852253, 1.02248405551571, 1.031289723437965, 1.021568528889323, 1.018752246233514, 1.011590410192247),
    @(1.02, 1.015420595870977, 1.027537545360704, 1.017125423170531, 1.012511024511812, 1.029988327489813, 1.021351816702886, 1.030727978147959, 1.020350827103264, 1.015752179693185, 1.011211701190615),
    @(1.02, 1.014108265968717, 1.02687234093594, 1.016020093298748, 1.010209940016578, 1.029800150449351, 1.020594940646718, 1.030350871721984, 1.019538630237741, 1.013750286963143, 1.010958302772797),
    @(1.02, 1.013540294137795, 1.026584610927721, 1.015542213860529, 1.009213626436927, 1.029717318537103, 1.020266722387712, 1.030186973479109, 1.019186852421299, 1.012882957347765, 1.010848360839859),
    @(1.02, 1.013329365921496, 1.026477783120709, 1.015364819542345, 1.008843556767267, 1.029686348669072, 1.020144734743144, 1.030126003881877, 1.019056173226375, 1.012560713524703, 1.010807490784012),
    @(1.02, 1.0133746088416, 1.026500695852324, 1.015402866171834, 1.008922937798984, 1.0296930009636, 1.020170904796207, 1.030139086147155, 1.019084204949701, 1.012629839562725, 1.010816259028499),
    @(1.02, 1.013522857890446, 1.0265757795205, 1.015527548113763, 1.009183036279064, 1.029714762684803, 1.020256640334074, 1.030181935553699, 1.019176050705995, 1.012856322209262, 1.010844983175122),
    @(1.02, 1.013614204574316, 1.026622047447434, 1.015604383557924, 1.009343292118986, 1.02972814399091, 1.020309455192961, 1.030208324524619, 1.019232638165322, 1.012995855046225, 1.010862676737333),
    @(1.02, 1.014145966230394, 1.026891443258889, 1.016051824156366, 1.01027606312902, 1.029805619322927, 1.020616713193858, 1.030361736366095, 1.019561974689856, 1.013807837934988, 1.010965594659838),
    @(1.02, 1.014479600353464, 1.027060511749794, 1.016332689392465, 1.010861181176688, 1.029853856409555, 1.020809318090706, 1.030457805277378, 1.019768534691441, 1.01431703726862, 1.011030093886618),
    @(1.02, 1.014674230065411, 1.027159156136228, 1.016496584259385, 1.011202477849022, 1.029881862001643, 1.020921614367378, 1.030513781856274, 1.019889008814707, 1.014613996977546, 1.011067694079447),
    @(1.02, 1.014740598262664, 1.027192796291828, 1.016552480177889, 1.011318852574282, 1.029891389073633, 1.020959896533239, 1.030532858421926, 1.019930085865681, 1.014715244614174, 1.011080511187866),
    @(1.02, 1.014443801818079, 1.027042369227425, 1.016302547868217, 1.01079840283802, 1.029848694499579, 1.020788658257481, 1.030447504066585, 1.019746373668225, 1.014262409958391, 1.011023175912997),
    @(1.02, 1.013479201083247, 1.026553667915181, 1.015490829307948, 1.009106443647489, 1.029708359984412, 1.020231395361982, 1.03016931996186, 1.019149004782139, 1.012789630930101, 1.010836525532198),
    @(1.02, 1.01287296041271, 1.026246680127249, 1.014981114476163, 1.008042671332159, 1.029618955241414, 1.019880601155658, 1.029993891482196, 1.018773338808615, 1.011863177154027, 1.010718981749224),
    @(1.02, 1.013194316984083, 1.026409393243774, 1.015251262545567, 1.008606596080921, 1.029666461232894, 1.020066603620817, 1.030086938704433, 1.018972493545192, 1.012354352804369, 1.010781311828813),
    @(1.02, 1.014459977553133, 1.027050566961111, 1.01631616730649, 1.010826769643494, 1.029851027347094, 1.02079799369169, 1.030452158922098, 1.019756387307981, 1.01428709383465, 1.011026301912729),
    @(1.02, 1.015929747062111, 1.027795768254838, 1.017554709134108, 1.013403458605117, 1.030060093715263, 1.021644890380276, 1.030873666878264, 1.020665704244386, 1.01652807838217, 1.01130977044374)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $r = $startRow + $i
    $vals = $rowData[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = "$($cols[$j])$r"
        $ws.Range($addr).Value = $vals[$j]
    }
}

Write-Output "Updated $($rowData.Length) rows"